$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.940.29'
$ws.Range('D3').Value = '2.389.81'
$ws.Range('E3').Value = '  -3.58%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.12'
$ws.Range('E5').Value = '  -1.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '88.78'
$ws.Range('E6').Value = '  -4.96%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.533'
$ws.Range('E7').Value = '  -3.62%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -4.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0842'
$ws.Range('E10').Value = '  -4.39%  '
$ws.Range('E11').Value = '  -6.20%  '
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').Value = '2.758.91'
$ws.Range('E13').Value = '  -3.55%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.63'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.30'
$ws.Range('E15').Value = '  -2.88%  '
$ws.Range('D16').Value = '2.425.94'
$ws.Range('E16').Value = '  -2.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.771'
$ws.Range('E17').Value = '  -3.55%  '
$ws.Range('D18').Value = '40.873.53'
$ws.Range('D19').Value = '0.0₃0918'
$ws.Range('E19').Value = '  -3.80%  '
$ws.Range('E20').Value = '  -4.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '69.29'
$ws.Range('E21').Value = '  -2.74%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.01'
$ws.Range('E22').Value = '  -2.88%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.88'
$ws.Range('E23').Value = '  -2.57%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.66'
$ws.Range('E24').Value = '  -4.08%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  -6.54%  '
$ws.Range('E27').Value = '  -3.07%  '
$ws.Range('E28').Value = '  -1.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.40'
$ws.Range('E29').Value = '  -4.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.19'
$ws.Range('E30').Value = '  -7.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '153.69'
$ws.Range('E31').Value = '  -2.39%  '
$ws.Range('E32').Value = '  +0.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.25'
$ws.Range('E33').Value = '  -4.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0738'
$ws.Range('E34').Value = '  -3.86%  '
$ws.Range('E35').Value = '  -5.10%  '
$ws.Range('E36').Value = '  -2.06%  '
$ws.Range('E37').Value = '  -4.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.14'
$ws.Range('E38').Value = '  -7.94%  '
$ws.Range('E39').Value = '  -4.18%  '
$ws.Range('E40').Value = '  -7.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.87'
$ws.Range('E41').Value = '  -3.93%  '
$ws.Range('E42').Value = '  -7.03%  '
$ws.Range('D43').Value = '1.978.79'
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0273'
$ws.Range('E44').Value = '  -4.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '17.69'
$ws.Range('E45').Value = '  -6.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.65'
$ws.Range('E46').Value = '  +1.63%  '
$ws.Range('E47').Value = '  -7.80%  '
$ws.Range('D48').Value = '2.619.97'
$ws.Range('E48').Value = '  -3.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '94.16'
$ws.Range('E49').Value = '  -4.37%  '
$ws.Range('E50').Value = '  -2.94%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.47'
$ws.Range('E51').Value = '  -2.24%  '
